{"js": "// 1) \"transfer\" -> \"transferIR\" in the checkCharacterType/quantLetter/countLetter bullet.\nconst searchResults = context.document.body.search(\n  \"checkStrengh e transfer para a classe StringUtils\",\n  { matchCase: true }\n);\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\n    \"checkStrengh e transferIR para a classe StringUtils\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 2) Add a new bulleted list item right after the \"getDigits\" paragraph,\n//    describing the rewrite of verifyIsCharEqual (formerly isCharEqual).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet getDigitsParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Cria\u00e7\u00e3o da fun\u00e7\u00e3o getDigits\") !== -1) {\n    getDigitsParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (getDigitsParagraph) {\n  // insertParagraph(\"After\") inherits the source paragraph's list/number\n  // formatting (numPr/ilvl/numId, indent, contextualSpacing, alignment,\n  // run properties), matching the existing list items in this section.\n  getDigitsParagraph.insertParagraph(\n    \"Reeecrita de verifyIsCharEqual (antiga isCharEqual), de modo a ficar mais leg\u00edvel e melhorar entendimento da l\u00f3gica da fun\u00e7\u00e3o\",\n    \"After\"\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"transfer\" -> \"transferIR\" in the checkCharacterType/quantLetter/countLetter bullet.\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Text = \"checkStrengh e transfer para a classe StringUtils\"\n$r.Find.Replacement.ClearFormatting()\n$r.Find.Replacement.Text = \"checkStrengh e transferIR para a classe StringUtils\"\n$r.Find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n\n# 2) Add a new bulleted list item right after the \"getDigits\" paragraph,\n#    describing the rewrite of verifyIsCharEqual (formerly isCharEqual).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Cria\u00e7\u00e3o da fun\u00e7\u00e3o getDigits*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    # InsertParagraphAfter() inherits the source paragraph's list/number\n    # formatting (numPr/ilvl/numId, indent, contextualSpacing, alignment,\n    # run properties), matching the existing list items in this section.\n    $null = $target.Range.InsertParagraphAfter()\n    $newPara = $target.Next()\n    $newPara.Range.Text = \"Reeecrita de verifyIsCharEqual (antiga isCharEqual), de modo a ficar mais leg\u00edvel e melhorar entendimento da l\u00f3gica da fun\u00e7\u00e3o\"\n}\n"}
